$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9359
$ws1.Range("F3").Value = 209
$ws1.Range("F5").Value = 511

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9359
$ws4.Range("F3").Value = 209
$ws4.Range("F5").Value = 511
